$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3831433893464293
$ws.Range("C2").Value = 0.05921439310084509
$ws.Range("D2").Value = 0.2021402031681134
$ws.Range("E2").Value = 0.1699234821611597
$ws.Range("F2").Value = 1.364870054439962
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("I2").Value = 0.6573968767895586
$ws.Range("J2").Value = 0.186067691903574
$ws.Range("K2").Value = 0.4153500648292834
$ws.Range("N2").Value = 1.534556796755314
$ws.Range("O2").Value = 3.250231393374804
$ws.Range("B3").Value = 0.346556099387044
$ws.Range("C3").Value = 0.05264474798234176
$ws.Range("D3").Value = 0.1959791953396888
$ws.Range("E3").Value = 0.1656344628850945
$ws.Range("F3").Value = 1.366129896015096
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("I3").Value = 0.6627194558780936
$ws.Range("J3").Value = 0.1820667567137662
$ws.Range("K3").Value = 0.3741560588014181
$ws.Range("N3").Value = 1.550230229227868
$ws.Range("O3").Value = 3.265215660862069
$ws.Range("B4").Value = 0.3241472087430566
$ws.Range("C4").Value = 0.04861178352787476
$ws.Range("D4").Value = 0.1922812381551466
$ws.Range("E4").Value = 0.1630856878434166
$ws.Range("F4").Value = 1.367601088904337
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("I4").Value = 0.6663089519006924
$ws.Range("J4").Value = 0.1797137774832862
$ws.Range("K4").Value = 0.3489117275349543
$ws.Range("N4").Value = 1.560356193173863
$ws.Range("O4").Value = 3.276085227156884
$ws.Range("B5").Value = 0.3150300246392135
$ws.Range("C5").Value = 0.04696859921641305
$ws.Range("D5").Value = 0.1907957639144655
$ws.Range("E5").Value = 0.1620683984033278
$ws.Range("F5").Value = 1.368376165317528
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("I5").Value = 0.6678525299136808
$ws.Range("J5").Value = 0.1787810027244845
$ws.Range("K5").Value = 0.3386373740623299
$ws.Range("N5").Value = 1.564608957767444
$ws.Range("O5").Value = 3.280934500691359
$ws.Range("B6").Value = 0.3135170231802533
$ws.Range("C6").Value = 0.04669576889911298
$ws.Range("D6").Value = 0.1905504024506257
$ws.Range("E6").Value = 0.1619007699043742
$ws.Range("F6").Value = 1.368515472839825
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("I6").Value = 0.6681137219911548
$ws.Range("J6").Value = 0.1786276924288259
$ws.Range("K6").Value = 0.3369321234714846
$ws.Range("N6").Value = 1.5653227579615
$ws.Range("O6").Value = 3.281765079315292
$ws.Range("B7").Value = 0.3240241912476733
$ws.Range("C7").Value = 0.04858962173145187
$ws.Range("D7").Value = 0.1922611174288278
$ws.Range("E7").Value = 0.1630718817704775
$ws.Range("F7").Value = 1.367610830882114
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("I7").Value = 0.6663294418489478
$ws.Range("J7").Value = 0.1797010921065834
$ws.Range("K7").Value = 0.348773110909292
$ws.Range("N7").Value = 1.560413035827921
$ws.Range("O7").Value = 3.276148926105691
$ws.Range("B8").Value = 0.37051683349614
$ws.Range("C8").Value = 0.05694905091365854
$ws.Range("D8").Value = 0.1999983296303753
$ws.Range("E8").Value = 0.1684270822683018
$ws.Range("F8").Value = 1.365159723253186
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("I8").Value = 0.659165397797068
$ws.Range("J8").Value = 0.1846666810161821
$ws.Range("K8").Value = 0.4011365619078333
$ws.Range("N8").Value = 1.539856731467598
$ws.Range("O8").Value = 3.255051596255157
$ws.Range("B9").Value = 0.4621111582371782
$ws.Range("C9").Value = 0.07334594469712385
$ws.Range("D9").Value = 0.2158405835517385
$ws.Range("E9").Value = 0.179598810505432
$ws.Range("F9").Value = 1.36588464142713
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("I9").Value = 0.6476667727051293
$ws.Range("J9").Value = 0.1952257254458658
$ws.Range("K9").Value = 0.5041879994405747
$ws.Range("N9").Value = 1.50353100970562
$ws.Range("O9").Value = 3.226921202591569
$ws.Range("B10").Value = 0.5296417641420703
$ws.Range("C10").Value = 0.08539309955045837
$ws.Range("D10").Value = 0.2278837548161476
$ws.Range("E10").Value = 0.1882137977710698
$ws.Range("F10").Value = 1.369786322194969
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("I10").Value = 0.6407732323686659
$ws.Range("J10").Value = 0.2034846020252132
$ws.Range("K10").Value = 0.5801020371219181
$ws.Range("N10").Value = 1.479269346664543
$ws.Range("O10").Value = 3.214325847330088
$ws.Range("B11").Value = 0.5604103722748448
$ws.Range("C11").Value = 0.09087335768751359
$ws.Range("D11").Value = 0.2334493591923064
$ws.Range("E11").Value = 0.1922211330145558
$ws.Range("F11").Value = 1.372292362265952
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("I11").Value = 0.6379747378997713
$ws.Range("J11").Value = 0.2073507274521376
$ws.Range("K11").Value = 0.6146770701683977
$ws.Range("N11").Value = 1.468758513699685
$ws.Range("O11").Value = 3.210348869711225
$ws.Range("B12").Value = 0.5720681318272511
$ws.Range("C12").Value = 0.09294852142045329
$ws.Range("D12").Value = 0.2355693324524424
$ws.Range("E12").Value = 0.1937512639851064
$ws.Range("F12").Value = 1.373346389085185
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("I12").Value = 0.6369635425223805
$ws.Range("J12").Value = 0.2088304059642212
$ws.Range("K12").Value = 0.6277751624267864
$ws.Range("N12").Value = 1.464853953951714
$ws.Range("O12").Value = 3.209094898844285
$ws.Range("B13").Value = 0.5695571521418117
$ws.Range("C13").Value = 0.09250160316241818
$ws.Range("D13").Value = 0.2351122088101789
$ws.Range("E13").Value = 0.1934211617280255
$ws.Range("F13").Value = 1.373114715202178
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("I13").Value = 0.6371791627590575
$ws.Range("J13").Value = 0.208511034702255
$ws.Range("K13").Value = 0.6249540298429679
$ws.Range("N13").Value = 1.465691504918897
$ws.Range("O13").Value = 3.209353755284212
$ws.Range("B14").Value = 0.5613693394224981
$ws.Range("C14").Value = 0.09104408496017413
$ws.Range("D14").Value = 0.2336235228827803
$ws.Range("E14").Value = 0.1923467648424051
$ws.Range("F14").Value = 1.372376972992711
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("I14").Value = 0.6378905734765397
$ws.Range("J14").Value = 0.2074721478071808
$ws.Range("K14").Value = 0.6157545561658253
$ws.Range("N14").Value = 1.468435766561035
$ws.Range("O14").Value = 3.21024065391228
$ws.Range("B15").Value = 0.5563548765942699
$ws.Range("C15").Value = 0.09015129863735183
$ws.Range("D15").Value = 0.2327132708265509
$ws.Range("E15").Value = 0.1916903099842173
$ws.Range("F15").Value = 1.371938761593682
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("I15").Value = 0.6383326540356471
$ws.Range("J15").Value = 0.2068378381538878
$ws.Range("K15").Value = 0.6101202803369858
$ws.Range("N15").Value = 1.470126560991559
$ws.Range("O15").Value = 3.210816725511194
$ws.Range("B16").Value = 0.5276318870486421
$ws.Range("C16").Value = 0.08503494359366925
$ws.Range("D16").Value = 0.2275217715028219
$ws.Range("E16").Value = 0.1879536811858245
$ws.Range("F16").Value = 1.36963725108059
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("I16").Value = 0.6409629113745865
$ws.Range("J16").Value = 0.2032341348127034
$ws.Range("K16").Value = 0.577843252385037
$ws.Range("N16").Value = 1.479966840348698
$ws.Range("O16").Value = 3.214621020069927
$ws.Range("B17").Value = 0.5100232599679657
$ws.Range("C17").Value = 0.08189615255064098
$ws.Range("D17").Value = 0.2243591750662972
$ws.Range("E17").Value = 0.185683960796851
$ws.Range("F17").Value = 1.368412546196566
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("I17").Value = 0.6426629184693482
$ws.Range("J17").Value = 0.2010513063255104
$ws.Range("K17").Value = 0.55805245742377
$ws.Range("N17").Value = 1.486138259767516
$ws.Range("O17").Value = 3.217403733728872
$ws.Range("B18").Value = 0.4998998553977003
$ws.Range("C18").Value = 0.08009080296254467
$ws.Range("D18").Value = 0.2225483395656624
$ws.Range("E18").Value = 0.1843867958557226
$ws.Range("F18").Value = 1.367776952909907
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("I18").Value = 0.6436724754996277
$ws.Range("J18").Value = 0.1998060739735763
$ws.Range("K18").Value = 0.5466732483197063
$ws.Range("N18").Value = 1.489737396549733
$ws.Range("O18").Value = 3.219169251472636
$ws.Range("B19").Value = 0.4964730559458417
$ws.Range("C19").Value = 0.07947954515535116
$ws.Range("D19").Value = 0.2219366350556697
$ws.Range("E19").Value = 0.1839490286439656
$ws.Range("F19").Value = 1.367573575488976
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("I19").Value = 0.6440197481798648
$ws.Range("J19").Value = 0.1993862249974114
$ws.Range("K19").Value = 0.5428211394767004
$ws.Range("N19").Value = 1.49096450290655
$ws.Range("O19").Value = 3.219795360381255
$ws.Range("B20").Value = 0.511897255638786
$ws.Range("C20").Value = 0.08223028281921074
$ws.Range("D20").Value = 0.2246949904443341
$ws.Range("E20").Value = 0.1859247159553234
$ws.Range("F20").Value = 1.368535795716113
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("I20").Value = 0.6424786628826133
$ws.Range("J20").Value = 0.2012826091638971
$ws.Range("K20").Value = 0.560158819159625
$ws.Range("N20").Value = 1.485476178347508
$ws.Range("O20").Value = 3.217090434967503
$ws.Range("B21").Value = 0.5637741310530942
$ws.Range("C21").Value = 0.09147219652447802
$ws.Range("D21").Value = 0.2340604505081956
$ws.Range("E21").Value = 0.1926619987015599
$ws.Range("F21").Value = 1.372590815655954
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("I21").Value = 0.637680297597381
$ws.Range("J21").Value = 0.2077768692914361
$ws.Range("K21").Value = 0.6184565266871118
$ws.Range("N21").Value = 1.467627656563657
$ws.Range("O21").Value = 3.209973310532718
$ws.Range("B22").Value = 0.5977154184745359
$ws.Range("C22").Value = 0.09751175168548798
$ws.Range("D22").Value = 0.2402535361436264
$ws.Range("E22").Value = 0.1971388493435455
$ws.Range("F22").Value = 1.375853237375367
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("I22").Value = 0.6348271869941691
$ws.Range("J22").Value = 0.2121125117478044
$ws.Range("K22").Value = 0.6565878997141397
$ws.Range("N22").Value = 1.456403666381473
$ws.Range("O22").Value = 3.206790841747448
$ws.Range("B23").Value = 0.5795971610617414
$ws.Range("C23").Value = 0.09428840715568754
$ws.Range("D23").Value = 0.2369416037177388
$ws.Range("E23").Value = 0.1947427518873752
$ws.Range("F23").Value = 1.374056033488785
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("I23").Value = 0.6363240558295686
$ws.Range("J23").Value = 0.2097901570260063
$ws.Range("K23").Value = 0.6362339033010187
$ws.Range("N23").Value = 1.462353751773639
$ws.Range("O23").Value = 3.208354977809051
$ws.Range("B24").Value = 0.5110500220557697
$ws.Range("C24").Value = 0.08207922506306886
$ws.Range("D24").Value = 0.2245431453049065
$ws.Range("E24").Value = 0.1858158464738082
$ws.Range("F24").Value = 1.368479861209806
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("I24").Value = 0.6425618645320483
$ws.Range("J24").Value = 0.2011780069197755
$ws.Range("K24").Value = 0.5592065367606835
$ws.Range("N24").Value = 1.48577534616779
$ws.Range("O24").Value = 3.217231561193302
$ws.Range("B25").Value = 0.4372894424898561
$ws.Range("C25").Value = 0.06890993377341204
$ws.Range("D25").Value = 0.2114835960724122
$ws.Range("E25").Value = 0.1765049615752403
$ws.Range("F25").Value = 1.365096661230687
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("I25").Value = 0.6505044641944693
$ws.Range("J25").Value = 0.1993862249974114
$ws.Range("K25").Value = 0.4762728017313123
$ws.Range("N25").Value = 1.470126560991559
$ws.Range("O25").Value = 3.210816725511194

Write-Output "Applied 380 kV case values"
